# Hortaliza / Coliflor - Feria Lagunitas de Puerto Montt
# Weekly update: insert two new observation rows (new week) and correct
# a mis-tagged Origen value on an existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert two new rows at position 136, pushing the existing rows
#    136:185 down to 138:187.
# ---------------------------------------------------------------------
$ws.Rows("136:137").Insert()

# ---------------------------------------------------------------------
# 2) Populate the two newly inserted rows with the new weekly data.
# ---------------------------------------------------------------------
$ws.Range("A136").Value = 4
$ws.Range("B136").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C136").Value = "Los Lagos"
$ws.Range("D136").Value = 44463
$ws.Range("E136").Value = 10
$ws.Range("F136").Value = 100112008
$ws.Range("G136").Value = "Coliflor"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 600
$ws.Range("K136").Value = 1300
$ws.Range("L136").Value = 1300
$ws.Range("M136").Value = 1300
$ws.Range("N136").Value = "`$/unidad"
$ws.Range("O136").Value = "Región Metropolitana"
$ws.Range("P136").Value = 1300
$ws.Range("Q136").Value = 1
$ws.Range("R136").Value = "Hortaliza"

$ws.Range("A137").Value = 4
$ws.Range("B137").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C137").Value = "Los Lagos"
$ws.Range("D137").Value = 44463
$ws.Range("E137").Value = 10
$ws.Range("F137").Value = 100112008
$ws.Range("G137").Value = "Coliflor"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Segunda"
$ws.Range("J137").Value = 600
$ws.Range("K137").Value = 1000
$ws.Range("L137").Value = 1000
$ws.Range("M137").Value = 1000
$ws.Range("N137").Value = "`$/unidad"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 1000
$ws.Range("Q137").Value = 1
$ws.Range("R137").Value = "Hortaliza"

# ---------------------------------------------------------------------
# 3) Data correction: the row that is now 149 (formerly row 147, with
#    Fecha = 44369 / Segunda... actually Primera / J=500) had the
#    wrong Origen; it should be "Región del Maule" instead of
#    "Región Metropolitana".
# ---------------------------------------------------------------------
$ws.Range("O149").Value = "Región del Maule"

# ---------------------------------------------------------------------
# 4) Refresh the sheet dimension to match the new data extent.
# ---------------------------------------------------------------------
$ws.UsedRange | Out-Null
